$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Wins"/"Looses" columns -------------------------------------------------

# Headers
$ws.Range("D1").Value = "Wins "
$ws.Range("E1").Value = "Looses"

# Style the existing header cells (A1:C1) with the new theme-color fill + white font
$hdrOld = $ws.Range("A1:C1")
$hdrOld.Font.ThemeColor = 2
$hdrOld.Interior.ThemeColor = 3

# Style the new header cells (D1:E1) the same way
$hdrNew = $ws.Range("D1:E1")
$hdrNew.Font.ThemeColor = 2
$hdrNew.Interior.ThemeColor = 3

# Data rows: everybody starts at 0 wins / 0 losses
$ws.Range("D2:D7").Value = 0
$ws.Range("E2:E7").Value = 0

# Mirror the author's final selection
[void]$ws.Range("D9").Select()
